# "formatted input for Yahoo add"
# - Add a new "Yahoo" sheet (after TTD) with segment data.
# - AppNexus: mark a segment's "Is Public" flag true and clear its Buyer Member ID.
# - TTD: remove the two rows whose segments are no longer offered
#   ("Healthy Food Buyers" and "Good Food"), and re-fit the Segment Full Path column.

$wb = $excel.ActiveWorkbook

$wsAppNexus = $wb.Worksheets.Item("AppNexus")
$wsTTD = $wb.Worksheets.Item("TTD")

# ---------------------------------------------------------------------------
# 1. AppNexus: row 3 (Eyeota .. Luggage) becomes publicly buyable, and loses
#    its Buyer Member ID value.
# ---------------------------------------------------------------------------
$wsAppNexus.Range("H3").Value = $true
$wsAppNexus.Range("J3").ClearContents()

# ---------------------------------------------------------------------------
# 2. TTD: drop the "Healthy Food Buyers" (row 4) and "Good Food" (row 6) rows.
#    Delete the lower row first so the earlier row index stays valid.
# ---------------------------------------------------------------------------
$wsTTD.Rows(6).Delete()
$wsTTD.Rows(4).Delete()

# Re-fit column H (Segment Full Path) now that the rows/content changed.
$wsTTD.Columns("H").AutoFit()

# ---------------------------------------------------------------------------
# 3. Add the new "Yahoo" worksheet after "TTD".
# ---------------------------------------------------------------------------
$wsYahoo = $wb.Worksheets.Add($null, $wsTTD)
$wsYahoo.Name = "Yahoo"

# Header row, copied in style from the TTD sheet's bold header cells.
$wsYahoo.Range("A1").Value = "Segment ID"
$wsYahoo.Range("B1").Value = "Segment Name"
$wsYahoo.Range("C1").Value = "Segment Description"
$wsYahoo.Range("A1:C1").Font.Bold = $true

# "Required" marker row, highlighted the same green as the other sheets.
$wsYahoo.Range("A2:C2").Value = "Required"
$wsYahoo.Range("A2:C2").Interior.Color = 5296274

# New Yahoo segment rows: Segment ID + Segment Description entered first,
# Segment Name filled in afterwards.
$wsYahoo.Range("A3").Value = 4444
$wsYahoo.Range("C3").Value = "UK Kantar Media TGI - Grocery Shopping - Premium Range Purchasers"

$wsYahoo.Range("A4").Value = 1111
$wsYahoo.Range("C4").Value = "Purchase Category - Reach - Healthy Food Buyers - Damn Rich People"

$wsYahoo.Range("A5").Value = 2222
$wsYahoo.Range("C5").Value = "Purchase Category - Reach - Healthy Food Buyers - Damn Poor People"

$wsYahoo.Range("A6").Value = 3333
$wsYahoo.Range("C6").Value = "Purchase Category - Reach - Healthy Food Buyers - Damn Poor People - Even Poorer People"

$wsYahoo.Range("B3").Value = "Just random segment 4444"
$wsYahoo.Range("B4").Value = "Just random segment 1111"
$wsYahoo.Range("B5").Value = "Just random segment 2222"
$wsYahoo.Range("B6").Value = "Just random segment 3333"

# ---------------------------------------------------------------------------
# 4. Restore per-sheet selections, then land on Yahoo as the active tab.
# ---------------------------------------------------------------------------
$wsAppNexus.Activate()
$wsAppNexus.Range("A3").Select()

$wsTTD.Activate()
$wsTTD.Range("H10").Select()

$wsYahoo.Activate()
$wsYahoo.Range("C7").Select()
